$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing Late / heading / Outstanding columns one place to the right.
# The newly inserted column inherits the width of the column to its left.
$inheritedWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $inheritedWidth

# Make "Repayment schedule" the active sheet/tab and set the new
# selection on it.
$ws.Select() | Out-Null
$ws.Range("S10").Select() | Out-Null
